$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("relays")

# Update relay settings values (demo of October 2016 update)
$ws.Range("E2").Value = 8
$ws.Range("F2").Value = 15
$ws.Range("G2").Value = 0.8

$ws.Range("D8").Value = 9000

$ws.Range("D9").Value = 10000

# Update the active selection on the sheet
$ws.Range("H2").Select()
